$p = $ppt.ActivePresentation
$s = $p.Slides.Item(28)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 2: "Notice the output from " / "your skill"
#     -> "Notice the output from your " / "skill"
$run1 = $tr.Find("Notice the output from ", 0)
$run1.InsertAfter("your ")

$run2 = $tr.Find("your skill", 0)
$run2.Delete()

$para2 = $tr.Paragraphs(2, 1)
$para2.InsertAfter("skill")

# --- New paragraph 3: "Do you have the Alexa app on your cell phone?"
$para2b = $tr.Paragraphs(2, 1)
$para2b.InsertAfter([char]13 + "Do you have the Alexa app on your cell phone?")

# --- New paragraph 4: "Output from your skill will appear there,too" (lvl 1)
$para3 = $tr.Paragraphs(3, 1)
$para3.InsertAfter([char]13 + "Output from your skill will appear there,too")

$para4 = $tr.Paragraphs(4, 1)
$para4.IndentLevel = 2

Write-Host "para4 text: [$($para4.Text)] len=$($para4.Text.Length)"

# "Output from your skill will appear " is 36 chars, "there,too" is 9 chars
# (para4.Text includes a trailing CR, so its printable length is 45)
$firstLen = "Output from your skill will appear ".Length
Write-Host "firstLen:" $firstLen

$secondRun = $para4.Characters($firstLen + 1, 9)
Write-Host "secondRun text: [$($secondRun.Text)]"
$secondRun.Delete()

$firstRun = $para4.Characters(1, $firstLen)
Write-Host "firstRun text (after delete): [$($firstRun.Text)]"
$firstRun.InsertAfter("there,too")

$para4 = $tr.Paragraphs(4, 1)
Write-Host "para4 final text: [$($para4.Text)]"

$firstRunFinal = $para4.Characters(1, $firstLen)
$firstRunFinal.Font.Size = 30.5

$secondRunFinal = $para4.Characters($firstLen + 1, 9)
Write-Host "secondRunFinal text: [$($secondRunFinal.Text)]"
$secondRunFinal.Font.Size = 30.5
